$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "VENTAS POR GRUPO"
# A new client ("ABREU FERNANDEZ JOSE PABLO") is inserted immediately
# before the existing "AGUILERA ANDRADE FAUSTO ROGELIO" row (row 282),
# pushing that whole block (and the trailing summary row) down by one.
# Separately, a monthly figure for an existing client (row 5, LAVABOS
# column I) is updated.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Direct value edit: I5 0 -> 158.4 (unrelated to the row insertion below)
$ws1.Cells.Item(5, 9).Value = 158.4

# Insert a new blank row at 282 (shifts 282..364 down to 283..365)
$ws1.Rows.Item(282).Insert()

# Populate the newly inserted row 282
$ws1.Cells.Item(282, 1).Value = "OFICINA-CATAECSA"
$ws1.Cells.Item(282, 2).Value = "ABREU FERNANDEZ JOSE PABLO"
for ($col = 3; $col -le 18; $col++) {
    $ws1.Cells.Item(282, $col).Value = 0
}

# The trailing "count" summary row (previously row 364, now row 365)
# reports "<count> de 362" per column; bump the denominator to 363 for
# every column, and bump column I's numerator from 14 to 15 (the new
# non-zero I5 value above grows that count).
$ws1.Cells.Item(365, 3).Value = "3 de 363"
$ws1.Cells.Item(365, 4).Value = "18 de 363"
$ws1.Cells.Item(365, 5).Value = "6 de 363"
$ws1.Cells.Item(365, 6).Value = "0 de 363"
$ws1.Cells.Item(365, 7).Value = "0 de 363"
$ws1.Cells.Item(365, 8).Value = "7 de 363"
$ws1.Cells.Item(365, 9).Value = "15 de 363"
$ws1.Cells.Item(365, 10).Value = "1 de 363"
$ws1.Cells.Item(365, 11).Value = "8 de 363"
$ws1.Cells.Item(365, 12).Value = "13 de 363"
$ws1.Cells.Item(365, 13).Value = "34 de 363"
$ws1.Cells.Item(365, 14).Value = "1 de 363"
$ws1.Cells.Item(365, 15).Value = "4 de 363"
$ws1.Cells.Item(365, 16).Value = "2 de 363"
$ws1.Cells.Item(365, 17).Value = "1 de 363"
$ws1.Cells.Item(365, 18).Value = "0 de 363"

# ---------------------------------------------------------------------
# Sheet 2: "VENTA MENSUAL"
# Same new client inserted before "AGUILERA ANDRADE FAUSTO ROGELIO"
# (row 286 here), and the same client's December ("diciembre") figure
# on row 5 is updated.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

# Direct value edit: F5 2218.33 -> 2376.73
$ws2.Cells.Item(5, 6).Value = 2376.73

# Insert a new blank row at 286 (shifts 286..368 down to 287..369)
$ws2.Rows.Item(286).Insert()

# Populate the newly inserted row 286
$ws2.Cells.Item(286, 1).Value = "OFICINA-CATAECSA"
$ws2.Cells.Item(286, 2).Value = "ABREU FERNANDEZ JOSE PABLO"
for ($col = 3; $col -le 7; $col++) {
    $ws2.Cells.Item(286, $col).Value = 0
}

# The trailing grand-total row (previously row 368, now row 369) holds
# literal totals; only the "diciembre" (F) total moves, by the same
# +158.4 delta as the F5 edit above.
$ws2.Cells.Item(369, 3).Value = 405958.18
$ws2.Cells.Item(369, 4).Value = 424433.47
$ws2.Cells.Item(369, 5).Value = 412473.7
$ws2.Cells.Item(369, 6).Value = 126324.97
$ws2.Cells.Item(369, 7).Value = 373790
